# Replace the real investor names in column A (A2:A6) with generic
# placeholder names ("Investor 1" .. "Investor 5") used for fund-upload
# tests, per commit message "Added tests for fund uploads".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

$investors = @("Investor 1", "Investor 2", "Investor 3", "Investor 4", "Investor 5")
for ($i = 0; $i -lt $investors.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $investors[$i]
}

# Give the updated investor cells their own style (distinguishing them
# visually from the rest of the sheet) and move the active selection to
# the investor column that was just edited.
$ws.Range("A2:A6").Font.Size = 10
$ws.Range("A2:A6").Select() | Out-Null
